$d = $word.ActiveDocument

# --- 1. Insert the "Regular meeting time changed..." bullet right before
#        the "New work assignment." bullet, reusing the same list
#        paragraph formatting (ListParagraph style, numId 17). ---
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "New work assignment\.") {
        $target = $p
        break
    }
}

$newRange = $target.Range
$newRange.InsertParagraphBefore()

# Re-locate the freshly inserted (still empty) paragraph.
$target2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "New work assignment\.") {
        $target2 = $p
        break
    }
}
$newPara = $target2.Previous()
$newParaRange = $newPara.Range
$newParaRange.Text = "Regular meeting time changed from 9 :00 to 2:00 pm due to time conflict with class. "
$newParaRange.Font.NameBi = "Times New Roman (Body CS)"

# --- 2. Mark both inline pictures' runs as NoProofing (<w:noProof/>). ---
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shape = $d.InlineShapes.Item($i)
    $shape.Range.NoProofing = -1
}
